$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")
$ws.Activate() | Out-Null

# --- Update existing parameter values -------------------------------------------------

# End Year: 2050 -> 2025
$ws.Range("B4").Value = 2025

# maximum_investment_capacity_per_year: 1000000 -> 1000
$ws.Range("B13").Value = 1000

# realistic_candidate_capacities_for_future: stays TRUE, description reworded
$ws.Range("C17").Value = "If this is true, the real capacity of the candidate power plants is considered for the FUTURE testing. Otherwise the dummy capacity inidated "

# realistic_candidate_capacities_tobe_installed: TRUE -> FALSE, and clarify its description
$ws.Range("B16").Value = $false
$ws.Range("C16").Value = "If this is true, the real capacity  of the power plants is chosen to be installed"

# dummy_capacity: 100 -> 1000
$ws.Range("B19").Value = 1000

# targetinvestment_per_year: TRUE -> FALSE
$ws.Range("B20").Value = $false

# --- Add a small spacer row then two warning/checker rows ------------------------------

$ws.Rows.Item(22).RowHeight = 10

$ws.Range("B5").Copy()
$ws.Range("B24").PasteSpecial(-4122)
$ws.Range("B24").Formula = '=IF(AND(B17=FALSE,B16=TRUE),"DANGER","ok")'
$ws.Range("C24").Value = "Testing different capacity than the one being installed can cause deviations in reality"

$ws.Range("B5").Copy()
$ws.Range("B23").PasteSpecial(-4122)
$ws.Range("B23").Formula = '=IF(AND(B17=TRUE,B16=FALSE),"DANGER!!!!!","ok")'
$ws.Range("C23").Value = "If the dummy capacity will be installed, it could be very different than expected "

# --- Reset the view: scroll back to top, select B4 --------------------------------------

$ws.Range("B4").Select() | Out-Null
